$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Kazhakuttom"
$ws.Range("E2").Value = 4512457889
$ws.Range("G2").Value = "Student"

$ws.Range("E2").Select()
